$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update User Name (C3) and Password (C4) on the TEST DATA sheet
$ws1.Range("C3").Value = "karthik18"
$ws1.Range("C4").Value = "LXCW11"

# Reflect the new cell selection recorded for Sheet1
$ws1.Activate()
$ws1.Range("E4").Select()
